$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.357.24'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.877.95'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7212'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.66'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08017'
$ws.Range("E8").Value = '  +2.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3152'
$ws.Range("E9").Value = '  +1.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.00'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.877.23'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.75'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.232'
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7132'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("E16").Value = '  +5.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008502'
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.351.13'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.47'
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.761'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1593'
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.044'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.49'
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.501'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.413'
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.202'
$ws.Range("E31").Value = '  -7.32%  '
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.932'
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7627'
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.178'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.709'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01875'
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.277.91'
$ws.Range("E38").Value = '  +3.43%  '
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.445'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '113.00'
$ws.Range("E41").Value = '  +4.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9108'
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '74.20'
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("E44").Value = '  +6.83%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.024.63'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5228'
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.797'
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.510'
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4348'
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.097'
$ws.Range("E51").Value = '  +0.42%  '

Write-Host "Applied changes"
